# The recorded takeoff/landing timestamps in the first 10 per-session sheets
# were re-dated 3 days later (2020-05-04 -> 2020-05-07, same time-of-day),
# matching the commit "working version; next step group data by player".
# Row 1 (header) and unrelated rows/sheets are left untouched; only the
# existing date cells in row 2 (and row 3 for the two sheets that have a
# second flight) move forward by exactly 3 days.

$wb = $excel.ActiveWorkbook

# Sheet index (1-based, matches Sheet1..Sheet10) -> list of A1 cell refs
# holding a date/time serial that needs to shift by +3 days.
$targets = @{
    1  = @("A2", "B2")
    2  = @("A2", "A3")
    3  = @("A2", "B2")
    4  = @("A2")
    5  = @("A2", "A3")
    6  = @("A2", "B2")
    7  = @("A2", "B2")
    8  = @("A2", "B2")
    9  = @("A2", "B2")
    10 = @("A2", "B2")
}

foreach ($sheetIndex in $targets.Keys) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    foreach ($cellRef in $targets[$sheetIndex]) {
        $cell = $ws.Range($cellRef)
        $oldSerial = $cell.Value2()
        $cell.Value = $oldSerial + 3
    }
}
